$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.525.68'
$ws.Range('E2').Value = '  +4.90%  '
$ws.Range('D3').Value = '1.601.85'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.84'
$ws.Range('E5').Value = '  +2.16%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.00'
$ws.Range('E8').Value = '  +8.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.251'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').Value = '1.831.14'
$ws.Range('E12').Value = '  +2.67%  '
$ws.Range('D13').Value = '1.604.99'
$ws.Range('E13').Value = '  +3.07%  '
$ws.Range('E14').Value = '  +3.42%  '
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '28.545.69'
$ws.Range('E16').Value = '  +5.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.36'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.03'
$ws.Range('E18').Value = '  +7.20%  '
$ws.Range('E19').Value = '  +1.33%  '
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.96'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.60'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  +2.03%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('D34').Value = '1.420.72'
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('E36').Value = '  -4.46%  '
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('E39').Value = '  +2.83%  '
$ws.Range('E40').Value = '  +7.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.822'
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.75'
$ws.Range('E42').Value = '  -2.44%  '
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.984'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('E45').Value = '  +6.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.95'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').Value = '1.741.40'
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.52'
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('D50').Value = '0.0₆0108'
$ws.Range('E50').Value = '  +7.21%  '
$ws.Range('E51').Value = '  +0.58%  '
